$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.631.15'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.395.91'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.93%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.49%  '
$ws.Range("D9").Value = '2.400.94'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.343'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '2.828.83'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -1.44%  '
$ws.Range("D17").Value = '60.473.71'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '2.400.25'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '562.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.61%  '
$ws.Range("D29").Value = '2.513.65'
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").Value = '0.0₃0933'
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.14%  '
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  +2.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.369'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.22%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.82%  '
$ws.Range("D46").Value = '0.0₆0282'
$ws.Range("E46").Value = '  -1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.587'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.47%  '
